# Updated cryptos list on Tue Feb 28 03:40:54 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "23.460.80"

# Row 3
$ws.Range("D3").Value = "1.633.14"
$ws.Range("E3").Value = "  -0.51%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.000"
$ws.Range("E5").Value = "  +0.06%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.40"
$ws.Range("E6").Value = "  -0.94%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3745"
$ws.Range("E7").Value = "  -0.45%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3671"
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.77"
$ws.Range("E9").Value = "  -1.82%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08205"
$ws.Range("E10").Value = "  +0.01%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.231"
$ws.Range("E11").Value = "  -3.94%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9998"
$ws.Range("E12").Value = "  +0.03%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.60"
$ws.Range("E13").Value = "  -1.99%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.582"
$ws.Range("E14").Value = "  -1.52%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001252"
$ws.Range("E15").Value = "  -2.67%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.285"
$ws.Range("E16").Value = "  -2.05%  "

# Row 17
$ws.Range("D17").Value = "1.635.48"
$ws.Range("E17").Value = "  -0.47%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.42"
$ws.Range("E18").Value = "  -0.71%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06978"
$ws.Range("E19").Value = "  +0.84%  "

# Row 20
$ws.Range("E20").Value = "  -2.97%  "

# Row 21
$ws.Range("E21").Value = "  -1.65%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.0000"
$ws.Range("E22").Value = "  +0.21%  "

# Row 23
$ws.Range("E23").Value = "  -0.83%  "

# Row 24
$ws.Range("D24").Value = "23.474.69"
$ws.Range("E24").Value = "  -0.39%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.189"
$ws.Range("E25").Value = "  +3.18%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.42"
$ws.Range("E27").Value = "  +0.38%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.01"
$ws.Range("E28").Value = "  -1.05%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.327"
$ws.Range("E29").Value = "  -0.56%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.82"
$ws.Range("E30").Value = "  -0.96%  "

# Row 31
$ws.Range("D31").Value = "1.814.59"
$ws.Range("E31").Value = "  -0.66%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.272"
$ws.Range("E32").Value = "  -4.80%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.838"
$ws.Range("E33").Value = "  +0.09%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.029"
$ws.Range("E34").Value = "  +4.84%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.96"
$ws.Range("E35").Value = "  +4.89%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02791"
$ws.Range("E36").Value = "  -2.06%  "

# Row 37
$ws.Range("E37").Value = "  -0.56%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.07160"
$ws.Range("E38").Value = "  -3.38%  "

# Row 39
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.080"
$ws.Range("E39").Value = "  -2.30%  "

# Row 40
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08761"
$ws.Range("E40").Value = "  -1.63%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7092"
$ws.Range("E41").Value = "  -0.80%  "

# Row 42
$ws.Range("E42").Value = "  -2.58%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.31"
$ws.Range("E43").Value = "  +0.28%  "

# Row 44
$ws.Range("E44").Value = "  -2.17%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6570"
$ws.Range("E45").Value = "  -0.14%  "

# Row 46
$ws.Range("E46").Value = "  -0.66%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9990"
$ws.Range("E47").Value = "  +0.07%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.002"
$ws.Range("E48").Value = "  -1.14%  "

# Row 49
$ws.Range("E49").Value = "  +0.60%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.212"
$ws.Range("E50").Value = "  +0.09%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "125.79"
$ws.Range("E51").Value = "  -3.50%  "
